# Added multi classes support:
# H1 previously held "班别"; replace it with the new column label so the
# template now instructs uploaders how to supply multiple classes.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H1").Value = "classes (class1,class2,class3)"

# Leave the selection on H1, matching the edited cell.
$ws.Range("H1").Select()
